$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 430.36

$ws.Range("F23").Value = 352
$ws.Range("G23").Value = 18032.96

$ws.Range("B32").Value = 66452
$ws.Range("F32").Value = 64
$ws.Range("G32").Value = 1967.36

$ws.Range("B33").Value = 51755
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 30.74

$ws.Range("B40").Value = 72006.95

$ws.Range("F42").Value = 135
$ws.Range("G42").Value = 26563.95

$ws.Range("F44").Value = 538
$ws.Range("G44").Value = 19588.58

$ws.Range("F47").Value = 250
$ws.Range("G47").Value = 48222.5

$ws.Range("F59").Value = 113
$ws.Range("G59").Value = 6667

$ws.Range("F67").Value = 3
$ws.Range("G67").Value = 56.13

$ws.Range("B73").Value = 258852.68

$ws.Range("F123").Value = 191
$ws.Range("G123").Value = 21426.38

$ws.Range("F139").Value = 11
$ws.Range("G139").Value = 1657.37

$ws.Range("B145").Value = 89482.67999999999

$ws.Range("F185").Value = 53
$ws.Range("G185").Value = 7068.08

$ws.Range("B189").Value = 44956.98

$ws.Range("F219").Value = 129
$ws.Range("G219").Value = 8194.08

$ws.Range("F223").Value = 42
$ws.Range("G223").Value = 12180.42

$ws.Range("F236").Value = 61
$ws.Range("G236").Value = 2640.08

$ws.Range("F245").Value = 36
$ws.Range("G245").Value = 2156.04

$ws.Range("B247").Value = 89665.77

$ws.Range("F249").Value = 1
$ws.Range("G249").Value = 31.03

$ws.Range("B252").Value = 1186.99

$ws.Range("F263").Value = 104
$ws.Range("G263").Value = 6739.2

$ws.Range("B270").Value = 7491.67

$ws.Range("F273").Value = 20
$ws.Range("G273").Value = 2040.2

$ws.Range("F274").Value = 101
$ws.Range("G274").Value = 8998.09

$ws.Range("B280").Value = 102985.73

$ws.Range("F284").Value = 1816
$ws.Range("G284").Value = 33596

$ws.Range("F289").Value = 85
$ws.Range("G289").Value = 9741

$ws.Range("B291").Value = 52784.81

$ws.Range("F322").Value = 35
$ws.Range("G322").Value = 11053

$ws.Range("F329").Value = 93
$ws.Range("G329").Value = 7891.98

$ws.Range("F332").Value = 26
$ws.Range("G332").Value = 2696.2

$ws.Range("F336").Value = 160
$ws.Range("G336").Value = 6796.8

$ws.Range("F351").Value = 23
$ws.Range("G351").Value = 1159.43

$ws.Range("F358").Value = 18
$ws.Range("G358").Value = 3015

$ws.Range("F365").Value = 19
$ws.Range("G365").Value = 1922.8

$ws.Range("B367").Value = 64983
$ws.Range("C367").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F367").Value = 6
$ws.Range("G367").Value = 514.08

$ws.Range("B368").Value = 66194
$ws.Range("C368").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F368").Value = 35
$ws.Range("G368").Value = 2998.8

$ws.Range("B369").Value = 66196
$ws.Range("C369").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F369").Value = 28
$ws.Range("G369").Value = 2455.6

$ws.Range("B370").Value = 64985
$ws.Range("C370").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F370").Value = 13
$ws.Range("G370").Value = 1140.1

$ws.Range("B372").Value = 142629.3

$ws.Range("B375").Value = 61610
$ws.Range("E375").Value = 122.71
$ws.Range("F375").Value = -58
$ws.Range("G375").Value = -5957.18

$ws.Range("B376").Value = 63565
$ws.Range("E376").Value = 109.19
$ws.Range("F376").Value = 60
$ws.Range("G376").Value = 6162.6

$ws.Range("B397").Value = 60325
$ws.Range("E397").Value = 151.57
$ws.Range("F397").Value = -102
$ws.Range("G397").Value = -12939.72

$ws.Range("B398").Value = 63560
$ws.Range("E398").Value = 134.87
$ws.Range("F398").Value = 1
$ws.Range("G398").Value = 126.86

$ws.Range("F405").Value = 155
$ws.Range("G405").Value = 26556.15

$ws.Range("B409").Value = 33094.07

$ws.Range("F446").Value = 232
$ws.Range("G446").Value = 33556.48

$ws.Range("B451").Value = 92567.37

$ws.Range("F456").Value = 342
$ws.Range("G456").Value = 48081.78

$ws.Range("B458").Value = 106981

$ws.Range("F466").Value = 71
$ws.Range("G466").Value = 2273.42

$ws.Range("B470").Value = 84554.44

$ws.Range("F500").Value = 132
$ws.Range("G500").Value = 23913.12

$ws.Range("F503").Value = 109
$ws.Range("G503").Value = 4418.86

$ws.Range("B519").Value = 197160.36

$ws.Range("B550").Value = 53263
$ws.Range("E550").Value = 15.29
$ws.Range("F550").Value = -313
$ws.Range("G550").Value = -4009.53

$ws.Range("B551").Value = 65066
$ws.Range("E551").Value = 13.61
$ws.Range("F551").Value = 0
$ws.Range("G551").Value = 0

$ws.Range("B569").Value = 65067
$ws.Range("E569").Value = 15.65
$ws.Range("F569").Value = 0
$ws.Range("G569").Value = 0

$ws.Range("B570").Value = 53595
$ws.Range("E570").Value = 17.61
$ws.Range("F570").Value = -338
$ws.Range("G570").Value = -4978.74

$ws.Range("F608").Value = 112
$ws.Range("G608").Value = 31691.52

$ws.Range("B612").Value = 128626.98

$ws.Range("F626").Value = 13
$ws.Range("G626").Value = 726.5700000000001

$ws.Range("B633").Value = 59307.86

$ws.Range("F736").Value = 304
$ws.Range("G736").Value = 37042.4

$ws.Range("F741").Value = 68
$ws.Range("G741").Value = 5465.84

$ws.Range("B743").Value = 44981.75

$ws.Range("F748").Value = 14
$ws.Range("G748").Value = 4570.44

$ws.Range("F750").Value = 11
$ws.Range("G750").Value = 2762.1

$ws.Range("B756").Value = 50615.29

$ws.Range("F782").Value = 6
$ws.Range("G782").Value = 15625.2

$ws.Range("B787").Value = 32880.27

$ws.Range("F814").Value = 7
$ws.Range("G814").Value = 313.95

$ws.Range("F818").Value = 52
$ws.Range("G818").Value = 7424.04

$ws.Range("F823").Value = 14
$ws.Range("G823").Value = 2115.82

$ws.Range("F824").Value = 12
$ws.Range("G824").Value = 396.72

$ws.Range("F827").Value = 38
$ws.Range("G827").Value = 14356.78

$ws.Range("F828").Value = 542
$ws.Range("G828").Value = 55777.22

$ws.Range("F832").Value = 392
$ws.Range("G832").Value = 14437.36

$ws.Range("F833").Value = 75
$ws.Range("G833").Value = 3540.75

$ws.Range("B839").Value = 278756.1

$ws.Range("F879").Value = 61
$ws.Range("G879").Value = 2181.97

$ws.Range("F881").Value = 8
$ws.Range("G881").Value = 356.96

$ws.Range("B884").Value = 20377.79

$ws.Range("F889").Value = 102
$ws.Range("G889").Value = 3083.46

$ws.Range("F890").Value = 1497
$ws.Range("G890").Value = 244175.67

$ws.Range("F892").Value = 54
$ws.Range("G892").Value = 15274.98

$ws.Range("F893").Value = 55
$ws.Range("G893").Value = 7955.75

$ws.Range("B896").Value = 271247.22

$ws.Range("F938").Value = 129
$ws.Range("G938").Value = 13354.08

$ws.Range("B940").Value = 20688.05

$ws.Range("B941").Value = 3958265.92

$ws.Range("B942").Value = 3958265.92
